$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.480.86'
$ws.Range('E2').Value = '  +4.42%  '
$ws.Range('D3').Value = '3.332.47'
$ws.Range('E3').Value = '  +4.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.71'
$ws.Range('E5').Value = '  +2.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.57'
$ws.Range('E6').Value = '  +4.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.118'
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.435'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '3.907.05'
$ws.Range('E12').Value = '  +4.60%  '
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.79'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = '62.574.87'
$ws.Range('E16').Value = '  +4.57%  '
$ws.Range('D17').Value = '3.335.66'
$ws.Range('E17').Value = '  +4.39%  '
$ws.Range('E18').Value = '  +4.90%  '
$ws.Range('E19').Value = '  +4.90%  '
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '385.68'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.539'
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.87'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '0.0₃0962'
$ws.Range('E27').Value = '  +6.26%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.98'
$ws.Range('E29').Value = '  +3.04%  '
$ws.Range('E30').Value = '  +4.18%  '
$ws.Range('E31').Value = '  +2.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('E32').Value = '  +2.06%  '
$ws.Range('E33').Value = '  +7.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.71'
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.28'
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('E36').Value = '  +9.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.87'
$ws.Range('E37').Value = '  +10.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.36'
$ws.Range('E38').Value = '  +6.56%  '
$ws.Range('D39').Value = '2.856.78'
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('E41').Value = '  +7.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.33'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.752'
$ws.Range('E43').Value = '  +3.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.52'
$ws.Range('E44').Value = '  +2.66%  '
$ws.Range('E45').Value = '  +3.23%  '
$ws.Range('E46').Value = '  +4.54%  '
$ws.Range('E47').Value = '  +6.52%  '
$ws.Range('E48').Value = '  +3.56%  '
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.804'
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '284.07'
$ws.Range('E51').Value = '  +8.84%  '
